# Apply updated "dSF" (column F) values to the relevant rows.
# These values were repulled from source data; only column F (dSF) changes,
# column E (dS0) and the rest stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'F2'  = -3
    'F6'  = 0
    'F9'  = 3
    'F11' = 2
    'F15' = -2
    'F21' = -1
    'F24' = 0
    'F30' = 0
    'F36' = 2
    'F39' = 0
    'F41' = -1
    'F45' = 1
    'F47' = 2
    'F50' = 0
    'F55' = 4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
